# Updated solution for Tutorial 6
# Replace "/" with "-" in the attendance date strings (column A) and
# refresh the attendance counting columns (D..H) for rows 3-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 7;  Date = "11-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 12; Date = "29-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 13; Date = "01-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 14; Date = "05-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 15; Date = "08-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 16; Date = "12-09-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($item in $data) {
    $r = $item.Row
    $cellA = $ws.Cells.Item($r, 1)
    # Force text format so Excel does not auto-convert the dd-mm-yyyy
    # looking string into a real date serial number.
    $cellA.NumberFormat = "@"
    $cellA.Value = $item.Date
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
